# Repull data, push all data, mean calculation
# Updates the "dSF" column (F) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    10 = -1
    11 = 4
    14 = 0
    15 = -2
    16 = -3
    17 = 2
    18 = 6
    19 = 6
    21 = -3
    24 = -5
    25 = 1
    26 = 5
    28 = -1
    34 = -2
    35 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
